$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 18 de Abril de 2020 a las 14:22"
$ws.Range("B17").Value = 31589
$ws.Range("C17").Value = 1140
$ws.Range("E17").Value = 27738
$ws.Range("G17").Value = 142
$ws.Range("H17").Value = 3601
$ws.Range("A23").Value = "Suecia"
$ws.Range("B23").Value = 13822
$ws.Range("C23").Value = 606
$ws.Range("D23").Value = 550
$ws.Range("E23").Value = 11761
$ws.Range("F23").Value = 1054
$ws.Range("G23").Value = 111
$ws.Range("H23").Value = 1511
$ws.Range("A24").Value = "Peru"
$ws.Range("B24").Value = 13489
$ws.Range("D24").Value = 6541
$ws.Range("E24").Value = 6648
$ws.Range("F24").Value = 137
$ws.Range("H24").Value = 300
$ws.Range("A33").Value = "Dinamarca"
$ws.Range("B33").Value = 7242
$ws.Range("C33").Value = 169
$ws.Range("D33").Value = 3847
$ws.Range("E33").Value = 3049
$ws.Range("F33").Value = 76
$ws.Range("G33").Value = 10
$ws.Range("H33").Value = 346
$ws.Range("A34").Value = "Arabia Saudita"
$ws.Range("B34").Value = 7142
$ws.Range("D34").Value = 1049
$ws.Range("E34").Value = 6006
$ws.Range("F34").Value = 74
$ws.Range("H34").Value = 87
$ws.Range("E35").Value = 6797
$ws.Range("G35").Value = 2
$ws.Range("H35").Value = 163
$ws.Range("B63").Value = 1832
$ws.Range("C63").Value = 18
$ws.Range("D63").Value = 615
$ws.Range("F63").Value = 27
$ws.Range("G63").Value = 3
$ws.Range("H63").Value = 39
$ws.Range("D70").Value = 171
$ws.Range("E70").Value = 1275
$ws.Range("B73").Value = 1317
$ws.Range("C73").Value = 13
$ws.Range("D73").Value = 190
$ws.Range("E73").Value = 1057
$ws.Range("F73").Value = 27
$ws.Range("G73").Value = 4
$ws.Range("H73").Value = 70
$ws.Range("A117").Value = "Kenia"
$ws.Range("B117").Value = 262
$ws.Range("C117").Value = 16
$ws.Range("D117").Value = 60
$ws.Range("E117").Value = 190
$ws.Range("F117").Value = 2
$ws.Range("G117").Value = 1
$ws.Range("H117").Value = 12
$ws.Range("A118").Value = "Sri Lanka"
$ws.Range("B118").Value = 248
$ws.Range("C118").Value = 4
$ws.Range("D118").Value = 77
$ws.Range("E118").Value = 164
$ws.Range("F118").Value = 1
$ws.Range("H118").Value = 7
$ws.Range("D141").Value = 64
$ws.Range("E141").Value = 32
$ws.Range("A151").Value = "San Martin (Parte Holandesa)"
$ws.Range("B151").Value = 64
$ws.Range("C151").Value = 7
$ws.Range("D151").Value = 12
$ws.Range("E151").Value = 43
$ws.Range("F151").Value = 6
$ws.Range("H151").Value = 9
$ws.Range("A152").Value = "Guyana"
$ws.Range("B152").Value = 63
$ws.Range("D152").Value = 9
$ws.Range("E152").Value = 48
$ws.Range("F152").Value = 4
$ws.Range("H152").Value = 6
$ws.Range("A153").Value = "Islas Caimanes"
$ws.Range("B153").Value = 61
$ws.Range("D153").Value = 7
$ws.Range("E153").Value = 53
$ws.Range("F153").Value = 3
$ws.Range("H153").Value = 1
$ws.Range("A195").Value = "Islas Turcas y Caicos"
$ws.Range("D195").Value = 0
$ws.Range("F195").Value = 0
$ws.Range("H195").Value = 1
$ws.Range("A196").Value = "Montserrat"
$ws.Range("D196").Value = 2
$ws.Range("E196").Value = 9
$ws.Range("F196").Value = 1
$ws.Range("H196").Value = 0
